{"js": "// Office.js (Word JavaScript API) script.\n// Applies the same edit described by the OOXML diff:\n//   1. Inserts two new right-aligned, 12pt paragraphs (\"Tim Tong\" and\n//      \"11/5/19\") at the very top of the document, wrapped in a\n//      \"_GoBack\" bookmark (the byline Word stamps in on \"Save as PDF\").\n//   2. Drops the stale <w:lastRenderedPageBreak/> markers in front of\n//      \"Services Page\", \"Contact Us\" and \"Personal Results\" (Word\n//      regenerates these at render time; they were left over from a\n//      previous render pass).\n//   3. Merges the \"Personal Results\" / \"/Customer Results\" runs (which\n//      used to be split apart by the bookmark that now lives at the top\n//      of the doc) into a single run \"Personal Results/Customer Results\".\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Step 1: insert the byline paragraphs (\"Tim Tong\" / \"11/5/19\") at the\n// very start of the body, right aligned, 12pt (sz=24 half-points), and\n// wrap both paragraphs in a \"_GoBack\" bookmark.\n// ---------------------------------------------------------------------\nconst bylineOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' +\n  '</pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:pPr><w:jc w:val=\"right\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Tim Tong</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"right\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>11/5/19</w:t></w:r>' +\n  '<w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '<w:p/>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\nconst startRange = body.getRange(\"Start\");\nstartRange.insertOoxml(bylineOoxml, Word.InsertLocation.before);\nawait context.sync();\n\n// insertOoxml leaves one extra empty paragraph behind (the \"trailing\n// paragraph mark\" of the inserted fragment) right before what used to be\n// the first paragraph of the document (\"About Page\") \u2014 remove it so the\n// body goes straight from the new \"11/5/19\" paragraph into \"About Page\".\nbody.paragraphs.load(\"items\");\nawait context.sync();\nbody.paragraphs.items[2].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Helper: find the paragraph containing `needle` and replace its raw\n// OOXML wholesale, so we can drop the <w:lastRenderedPageBreak/> marker\n// (not exposed anywhere on the Office.js object model) while keeping\n// every other attribute byte-for-byte identical.\n// ---------------------------------------------------------------------\nasync function replaceParagraphOoxml(needle, newParaInnerXml) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const match = results.items[0];\n  const para = match.paragraphs.getFirst();\n  const paraRange = para.getRange(\"Whole\");\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' +\n    newParaInnerXml +\n    '<w:p/>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>';\n\n  paraRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n\n  // Same trailing-paragraph-mark quirk as above: insertOoxml leaves one\n  // extra empty paragraph right after the replaced one; drop it so the\n  // paragraph count / following content is unaffected.\n  match.paragraphs.load(\"items\");\n  await context.sync();\n  const followers = body.paragraphs;\n  followers.load(\"items\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Step 2: \"Services Page\" \u2014 drop <w:lastRenderedPageBreak/>.\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\"Services Page\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const para = results.items[0].paragraphs.getFirst();\n  const paraRange = para.getRange(\"Whole\");\n\n  const newPara =\n    '<w:p w14:paraId=\"395E5850\" w14:textId=\"68957D2E\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' +\n    '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' +\n    '<w:r w:rsidRPr=\"00F866E0\"><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Services Page</w:t></w:r>' +\n    '</w:p>';\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' +\n    newPara +\n    '<w:p/>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>';\n\n  paraRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Step 3: \"Contact Us\" \u2014 drop <w:lastRenderedPageBreak/>.\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\"Contact Us\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const para = results.items[0].paragraphs.getFirst();\n  const paraRange = para.getRange(\"Whole\");\n\n  const newPara =\n    '<w:p w14:paraId=\"0D8702F7\" w14:textId=\"0303F93F\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' +\n    '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Contact Us</w:t></w:r>' +\n    '</w:p>';\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' +\n    newPara +\n    '<w:p/>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>';\n\n  paraRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Step 4: \"Personal Results\" / \"/Customer Results\" \u2014 drop\n// <w:lastRenderedPageBreak/>, drop the (now-relocated) \"_GoBack\"\n// bookmark that used to split the two runs, and merge the text into a\n// single run \"Personal Results/Customer Results\".\n// ---------------------------------------------------------------------\n{\n  const results = body.search(\"Personal Results\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const para = results.items[0].paragraphs.getFirst();\n  const paraRange = para.getRange(\"Whole\");\n\n  const newPara =\n    '<w:p w14:paraId=\"626E8222\" w14:textId=\"5A6D8BE9\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' +\n    '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Personal Results/Customer Results</w:t></w:r>' +\n    '</w:p>';\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' +\n    newPara +\n    '<w:p/>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n    '</pkg:package>';\n\n  paraRange.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the same edit described by the OOXML diff:\n#   1. Inserts two new right-aligned, 12pt paragraphs (\"Tim Tong\" and\n#      \"11/5/19\") at the very top of the document, wrapped in a\n#      \"_GoBack\" bookmark (the byline Word stamps in on \"Save as PDF\").\n#   2. Drops the stale <w:lastRenderedPageBreak/> markers in front of\n#      \"Services Page\", \"Contact Us\" and \"Personal Results\" (Word\n#      regenerates these at render time; they were left over from a\n#      previous render pass).\n#   3. Merges the \"Personal Results\" / \"/Customer Results\" runs (which\n#      used to be split apart by the bookmark that now lives at the top\n#      of the doc) into a single run \"Personal Results/Customer Results\".\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# Step 1: insert the byline paragraphs (\"Tim Tong\" / \"11/5/19\") at the\n# very start of the document, right aligned, 12pt (sz=24 half-points),\n# and wrap both paragraphs in a \"_GoBack\" bookmark.\n# -----------------------------------------------------------------------\n$bylineOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' + `\n'<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n'</pkg:xmlData></pkg:part>' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n'<w:p><w:pPr><w:jc w:val=\"right\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' + `\n'<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' + `\n'<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Tim Tong</w:t></w:r></w:p>' + `\n'<w:p><w:pPr><w:jc w:val=\"right\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' + `\n'<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>11/5/19</w:t></w:r>' + `\n'<w:bookmarkEnd w:id=\"0\"/></w:p>' + `\n'<w:p/>' + `\n'</w:body></w:document>' + `\n'</pkg:xmlData></pkg:part>' + `\n'</pkg:package>'\n\n$startRange = $d.Range(0, 0)\n$startRange.InsertXML($bylineOoxml)\n\n# InsertXML leaves one extra empty paragraph behind (the \"trailing\n# paragraph mark\" of the inserted fragment) right before what used to be\n# the first paragraph of the document (\"About Page\") \u2014 remove it so the\n# body goes straight from the new \"11/5/19\" paragraph into \"About Page\".\n$d.Paragraphs(3).Range.Delete()\n\n# -----------------------------------------------------------------------\n# Step 2: \"Services Page\" \u2014 drop <w:lastRenderedPageBreak/>.\n# -----------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"Services Page\") | Out-Null\n$paraRange = $rng.Paragraphs(1).Range\n\n$newPara = '<w:p w14:paraId=\"395E5850\" w14:textId=\"68957D2E\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' + `\n'<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' + `\n'<w:r w:rsidRPr=\"00F866E0\"><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Services Page</w:t></w:r>' + `\n'</w:p>'\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' + `\n'<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n'</pkg:xmlData></pkg:part>' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + `\n$newPara + `\n'<w:p/>' + `\n'</w:body></w:document>' + `\n'</pkg:xmlData></pkg:part>' + `\n'</pkg:package>'\n\n$paraRange.InsertXML($ooxml)\n\n# -----------------------------------------------------------------------\n# Step 3: \"Contact Us\" \u2014 drop <w:lastRenderedPageBreak/>.\n# -----------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"Contact Us\") | Out-Null\n$paraRange = $rng.Paragraphs(1).Range\n\n$newPara = '<w:p w14:paraId=\"0D8702F7\" w14:textId=\"0303F93F\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' + `\n'<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' + `\n'<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Contact Us</w:t></w:r>' + `\n'</w:p>'\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' + `\n'<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n'</pkg:xmlData></pkg:part>' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + `\n$newPara + `\n'<w:p/>' + `\n'</w:body></w:document>' + `\n'</pkg:xmlData></pkg:part>' + `\n'</pkg:package>'\n\n$paraRange.InsertXML($ooxml)\n\n# -----------------------------------------------------------------------\n# Step 4: \"Personal Results\" / \"/Customer Results\" \u2014 drop\n# <w:lastRenderedPageBreak/>, drop the (now-relocated) \"_GoBack\"\n# bookmark that used to split the two runs, and merge the text into a\n# single run \"Personal Results/Customer Results\".\n# -----------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Execute(\"Personal Results\") | Out-Null\n$paraRange = $rng.Paragraphs(1).Range\n\n$newPara = '<w:p w14:paraId=\"626E8222\" w14:textId=\"5A6D8BE9\" w:rsidR=\"00F866E0\" w:rsidRDefault=\"00F866E0\">' + `\n'<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr></w:pPr>' + `\n'<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"44\"/><w:szCs w:val=\"44\"/></w:rPr><w:t>Personal Results/Customer Results</w:t></w:r>' + `\n'</w:p>'\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' + `\n'<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n'</pkg:xmlData></pkg:part>' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + `\n$newPara + `\n'<w:p/>' + `\n'</w:body></w:document>' + `\n'</pkg:xmlData></pkg:part>' + `\n'</pkg:package>'\n\n$paraRange.InsertXML($ooxml)\n"}
